$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# 1. Workbook view window height
# ---------------------------------------------------------------
$wb.Windows.Item(1).Height = 15960

# ---------------------------------------------------------------
# 2. Insert new row 20: "Max Sum of a Pair With Equal Sum of Digits"
#    (plain Insert inherits formatting from row 19 above it)
# ---------------------------------------------------------------
$ws.Rows.Item(20).Insert()
$ws.Range("A20").Value = 2342
$ws.Range("B20").Value = "Max Sum of a Pair With Equal Sum of Digits"
$ws.Range("C20").Value = "Medium"
$ws.Range("D20").Value = "Arrays,hashmap"
$ws.Range("E20").Value = 45700

# ---------------------------------------------------------------
# 3. Append new "Sliding Window (Basics)" section at the bottom
#    Rows: 31 blank, 32 header, 33 blank, 34 data
# ---------------------------------------------------------------
$ws.Rows.Item(31).Insert()
$ws.Rows.Item(32).Insert()
$ws.Rows.Item(33).Insert()
$ws.Rows.Item(34).Insert()

# Row 31 & 33 must stay empty gap rows - clear any inherited formatting
$ws.Range("A31:E31").ClearFormats()
$ws.Range("A33:E33").ClearFormats()

# Row 32: section header (copy header look from B14)
$ws.Range("B14").Copy()
$ws.Range("B32").PasteSpecial(-4122)
$ws.Range("B32").Value = "Sliding Window (Basics)"
$ws.Rows.Item(32).RowHeight = 32

# Row 34: data row (copy look from row 16 which already has correct per-column styles)
$ws.Range("A16:E16").Copy()
$ws.Range("A34:E34").PasteSpecial(-4122)
$ws.Range("A34").Value = 3
$ws.Range("B34").Value = "Longest Substring Without Repeating Characters"
$ws.Range("C34").Value = "Medium"
$ws.Range("D34").Value = "String,Sliding Window,Two pointers"
$ws.Range("E34").Value = 45700

# ---------------------------------------------------------------
# 4. Selection / active cell
# ---------------------------------------------------------------
$ws.Range("B36").Select()

Write-Host "edit complete"
